$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "13.88" -> "13.90" (fix a percentage), keeping the number split
# across two runs ("13." and "90") the way the canonical edit recorded it.
# ---------------------------------------------------------------------------
$found = $d.Content.Duplicate
$found.Find.Execute("13.88") | Out-Null

$numStart = $found.Start
$numEnd = $found.End

# The split point is right after "13." (3 characters in).
$splitPoint = $numStart + 3
$suffixRange = $d.Range($splitPoint, $numEnd)

# Temporarily bookmark the "88" tail so the run boundary at $splitPoint is
# preserved once we overwrite its text with "90" - this keeps "13." and "90"
# as two distinct runs instead of silently recombining with their neighbours.
$marker = $d.Bookmarks.Add("TmpSplitMarker", $suffixRange)

$suffixRange = $d.Range($splitPoint, $numEnd)
$suffixRange.Text = "90"

$d.Bookmarks("TmpSplitMarker").Delete()

# ---------------------------------------------------------------------------
# Change 2: drop the stray grammar-check markers around "all of" and fold
# the sentence back into a single run.
# ---------------------------------------------------------------------------
$sentence = "Include the scrum master, and all of the members of the group (marking those who are present)."
$d.Content.Find.Execute($sentence, $false, $false, $false, $false, $false, $true, 1, $false, $sentence, 2) | Out-Null
